$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.885.94'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '2.505.54'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  +0.06%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '323.03'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.29%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '108.20'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -1.25%  '
$ws.Range('E7').Value = '  -0.70%  '
$ws.Range('E8').Value = '  +0.03%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.559'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +3.01%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '40.16'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('E11').Value = '  +5.33%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.0814'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -0.72%  '
$ws.Range('E13').Value = '  +0.57%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '7.17'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('D15').Value = '2.897.14'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').Value = '2.503.36'
$ws.Range('E16').Value = '  +0.18%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.848'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = '47.757.45'
$ws.Range('E18').Value = '  +0.94%  '
$ws.Range('E19').Value = '  +2.78%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '6.59'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').Value = '0.0₃0941'
$ws.Range('E21').Value = '  -0.86%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '2.75'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +4.58%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '70.91'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -0.23%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '247.33'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -1.09%  '
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('E26').Value = '  +0.01%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '25.78'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('B28').Value = 'Cosmos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '10.21'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +1.54%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '2.19'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -4.44%  '
$ws.Range('E30').Value = '  +4.10%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '35.02'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -2.79%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '49.80'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -0.92%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '19.98'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('E34').Value = '  -1.83%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('E36').Value = '  -1.61%  '
$ws.Range('E37').Value = '  -1.44%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '4.68'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -1.62%  '
$ws.Range('E39').Value = '  -1.09%  '
$ws.Range('E40').Value = '  -0.67%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '22.19'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +3.88%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '118.86'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -2.98%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.17'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -3.62%  '
$ws.Range('D45').Value = '1.995.44'
$ws.Range('E45').Value = '  -0.02%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '3.11'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +1.02%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '2.03'
$c.Style = "Normal"
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.81'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +0.85%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '9.10'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('E50').Value = '  -2.54%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '56.45'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +1.13%  '
